$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "p_seco [kN/m3]" (C) and "p_saturado [kN/m3]" (D) values for rows 3 and 4
$ws.Range("C3").Value = 17
$ws.Range("D3").Value = 20

$ws.Range("C4").Value = 19
$ws.Range("D4").Value = 22
